$d = $word.ActiveDocument

# The diff removes every <w:contextualSpacing w:val="0"/> element that Word
# had written into each paragraph's <w:pPr>. That property is not exposed
# anywhere on the scripting object model (ParagraphFormat has no
# ContextualSpacing member in this host), so we fall back to round-tripping
# each paragraph's own Open XML: read it back with Range.WordOpenXML, strip
# the <w:contextualSpacing/> element out of its <w:pPr>, and feed the
# trimmed markup back in with Range.InsertXML on that exact paragraph's
# range (which replaces only that paragraph's content, leaving the rest of
# the document untouched).

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $range = $para.Range

    $fullXml = $range.WordOpenXML
    if ($fullXml -notmatch "contextualSpacing") {
        continue
    }

    $bodyMatch = [regex]::Match($fullXml, "<w:body>(.*)</w:body>", "Singleline")
    if (-not $bodyMatch.Success) {
        continue
    }
    $bodyContent = $bodyMatch.Groups[1].Value

    # The paragraph addressed by $range is always emitted first inside
    # <w:body> by this host, so grab that leading <w:p>...</w:p> chunk.
    $paraMatch = [regex]::Match($bodyContent, "^(<w:p\b.*?</w:p>)", "Singleline")
    if (-not $paraMatch.Success) {
        continue
    }
    $paraXml = $paraMatch.Groups[1].Value

    if ($paraXml -notmatch "contextualSpacing") {
        continue
    }

    $newParaXml = [regex]::Replace($paraXml, "<w:contextualSpacing[^/]*/>", "")

    $wrapped = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' + '<w:body>' + $newParaXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $range.InsertXML($wrapped)
}

Write-Output "done"
